$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to match the new title
$ws.Name = "暫收款日餘額前後差異比較表"

# Move selection to D22 as recorded in the saved view state
$ws.Range("D22").Select()
